$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.626.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.58%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.78%  "

# Row 7
$ws.Range("E7").Value = "  +1.08%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +1.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "

# Row 11
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0847"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.468.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.998.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.98%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.973"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.683.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "

# Row 21
$ws.Range("E21").Value = "  +0.37%  "

# Row 22
$ws.Range("E22").Value = "  +1.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.20%  "

# Row 26
$ws.Range("E26").Value = "  -1.71%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("E31").Value = "  +3.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0464"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.81%  "

# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.06%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.80%  "

# Row 41
$ws.Range("E41").Value = "  -0.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.61%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.060.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.33%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0355"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.64%  "
